$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every changed row: write each target value as a string-literal
# formula (="text") so it evaluates to Text instead of Excel inferring
# Number/Percent from a plain numeric-looking string. Then immediately
# copy -> paste-special-values that rows contiguous block so the live
# formula collapses to a static value/shared-string, matching the
# original inlineStr cells type while leaving style (no "s" attr) and
# every other cell untouched.

$ws.Range("D2").Formula = '="302.86"'
$ws.Range("E2").Formula = '="1.25%"'
$r = $ws.Range("D2:E2")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D3").Formula = '="32.68"'
$ws.Range("E3").Formula = '="4.11%"'
$r = $ws.Range("D3:E3")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D4").Formula = '="4.969"'
$ws.Range("E4").Formula = '="-2.72%"'
$r = $ws.Range("D4:E4")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E5").Formula = '="-1.26%"'
$r = $ws.Range("E5")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D6").Formula = '="1.971"'
$ws.Range("E6").Formula = '="-14.75%"'
$r = $ws.Range("D6:E6")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D7").Formula = '="7.873"'
$ws.Range("E7").Formula = '="0.82%"'
$r = $ws.Range("D7:E7")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D8").Formula = '="3.799"'
$ws.Range("E8").Formula = '="-1.65%"'
$r = $ws.Range("D8:E8")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D9").Formula = '="0.9268"'
$ws.Range("E9").Formula = '="0.41%"'
$r = $ws.Range("D9:E9")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D10").Formula = '="0.1773"'
$ws.Range("E10").Formula = '="1.20%"'
$r = $ws.Range("D10:E10")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D11").Formula = '="0.07850"'
$ws.Range("E11").Formula = '="3.48%"'
$r = $ws.Range("D11:E11")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D12").Formula = '="0.08672"'
$ws.Range("E12").Formula = '="-7.24%"'
$r = $ws.Range("D12:E12")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E13").Formula = '="4.78%"'
$r = $ws.Range("E13")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D14").Formula = '="0.1004"'
$ws.Range("E14").Formula = '="0.02%"'
$r = $ws.Range("D14:E14")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D15").Formula = '="0.001530"'
$ws.Range("E15").Formula = '="1.19%"'
$r = $ws.Range("D15:E15")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D16").Formula = '="0.005780"'
$ws.Range("E16").Formula = '="-0.63%"'
$r = $ws.Range("D16:E16")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D17").Formula = '="3.462"'
$r = $ws.Range("D17")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E18").Formula = '="-4.95%"'
$r = $ws.Range("E18")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D19").Formula = '="0.3330"'
$ws.Range("E19").Formula = '="1.79%"'
$r = $ws.Range("D19:E19")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D20").Formula = '="0.1319"'
$ws.Range("E20").Formula = '="0.86%"'
$r = $ws.Range("D20:E20")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D21").Formula = '="4.305"'
$ws.Range("E21").Formula = '="7.25%"'
$r = $ws.Range("D21:E21")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E22").Formula = '="17.14%"'
$r = $ws.Range("E22")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D23").Formula = '="0.04577"'
$ws.Range("E23").Formula = '="-0.96%"'
$r = $ws.Range("D23:E23")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E24").Formula = '="-2.09%"'
$r = $ws.Range("E24")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D25").Formula = '="0.004436"'
$ws.Range("E25").Formula = '="-0.97%"'
$r = $ws.Range("D25:E25")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E26").Formula = '="0.05%"'
$r = $ws.Range("E26")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D39").Formula = '="0.01706"'
$ws.Range("E39").Formula = '="-2.05%"'
$r = $ws.Range("D39:E39")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D40").Formula = '="0.04762"'
$ws.Range("E40").Formula = '="3.04%"'
$r = $ws.Range("D40:E40")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D41").Formula = '="0.007477"'
$ws.Range("E41").Formula = '="7.13%"'
$r = $ws.Range("D41:E41")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D42").Formula = '="0.1356"'
$ws.Range("E42").Formula = '="-0.33%"'
$r = $ws.Range("D42:E42")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D43").Formula = '="0.002340"'
$ws.Range("E43").Formula = '="6.83%"'
$r = $ws.Range("D43:E43")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D44").Formula = '="0.01163"'
$ws.Range("E44").Formula = '="12.89%"'
$r = $ws.Range("D44:E44")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D45").Formula = '="0.00006253"'
$ws.Range("E45").Formula = '="-0.21%"'
$r = $ws.Range("D45:E45")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("E46").Formula = '="0.00%"'
$r = $ws.Range("E46")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D47").Formula = '="0.003104"'
$ws.Range("E47").Formula = '="-61.10%"'
$r = $ws.Range("D47:E47")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D48").Formula = '="0.8234"'
$ws.Range("E48").Formula = '="10.30%"'
$r = $ws.Range("D48:E48")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D49").Formula = '="0.00002100"'
$ws.Range("E49").Formula = '="0.00%"'
$r = $ws.Range("D49:E49")
$r.Copy()
$r.PasteSpecial(-4163)

$ws.Range("D50").Formula = '="0.0002000"'
$ws.Range("E50").Formula = '="0.00%"'
$r = $ws.Range("D50:E50")
$r.Copy()
$r.PasteSpecial(-4163)

$excel.CutCopyMode = $false
